# Log_of_all_Blogs.xlsx - append Post43 entry to the blog log table.
#
# Mirrors the authoring diff: a new row 53 is added to the bottom of the
# "Table2" table on Sheet1 (B10:F52 -> B10:F53), three new shared strings are
# introduced for the title + the two links, and the sheet view's selection
# moves to the newly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 52 (values + formatting) and insert it as the new row 53 so the
# new row picks up the same cell styles (S.No/Title/Date/links) as the rest
# of the table body.
$ws.Range("B52:F52").Copy()
$ws.Range("B53:F53").Insert()

# Overwrite the copied values with the new post's data.
$ws.Range("B53").Value = 43
$ws.Range("C53").Value = "Append Data in File | Shell Scripting "
$ws.Range("D53").Value = 44170
$ws.Range("E53").Value = "https://programmingport.hashnode.dev/append-data-in-file-or-shell-scripting"
$ws.Range("F53").Value = "https://dev.to/rahulmishra05/append-data-in-file-shell-scripting-ne0"

# Grow the table (Table2) so its range/autofilter covers the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B10:F53"))

# Match the author's final selection/scroll position on save.
$ws.Range("F53").Select()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 6
